$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new "June 2019" entry as the next row below "May 2019".
$ws.Range("A11").Value = "June 2019"
$ws.Range("B11").Value = "https://myemail.constantcontact.com/News-From-The-Forest---June-2019.html?soid=1102494320279&aid=qsCq9FpINss"

# Turn the URL cell into a real hyperlink, matching the style used by the
# other month rows (the Hyperlink cell style).
$ws.Hyperlinks.Add($ws.Range("B11"), "https://myemail.constantcontact.com/News-From-The-Forest---June-2019.html?soid=1102494320279&aid=qsCq9FpINss")
$ws.Range("B11").Style = "Hyperlink"

# Leave the selection where Excel would after this edit.
$ws.Range("B26").Select()
